# Update gh-pages to output generated at 456a3b4
# Regenerated "F" column (remaining ticket / stock count) values across the
# three data sheets (展览, 演出, 全部类型) to match the newly scraped data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 57
$ws1.Range("F3").Value  = 3288
$ws1.Range("F5").Value  = 2371
$ws1.Range("F7").Value  = 330
$ws1.Range("F8").Value  = 1345
$ws1.Range("F9").Value  = 1064
$ws1.Range("F11").Value = 499
$ws1.Range("F16").Value = 8292
$ws1.Range("F17").Value = 360
$ws1.Range("F18").Value = 2478
$ws1.Range("F21").Value = 171
$ws1.Range("F23").Value = 567
$ws1.Range("F26").Value = 993
$ws1.Range("F27").Value = 1920
$ws1.Range("F28").Value = 1534
$ws1.Range("F29").Value = 64
$ws1.Range("F31").Value = 1912
$ws1.Range("F34").Value = 20
$ws1.Range("F35").Value = 70
$ws1.Range("F38").Value = 52
$ws1.Range("F39").Value = 214
$ws1.Range("F40").Value = 389
$ws1.Range("F41").Value = 53
$ws1.Range("F43").Value = 245

# ---------------------------------------------------------------
# Sheet "演出" (Performances)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 19
$ws2.Range("F6").Value = 3

# ---------------------------------------------------------------
# Sheet "全部类型" (All types - combined view)
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 19
$ws4.Range("F3").Value  = 57
$ws4.Range("F4").Value  = 3288
$ws4.Range("F6").Value  = 2371
$ws4.Range("F8").Value  = 330
$ws4.Range("F9").Value  = 1345
$ws4.Range("F11").Value = 1064
$ws4.Range("F13").Value = 499
$ws4.Range("F17").Value = 8292
$ws4.Range("F18").Value = 360
$ws4.Range("F19").Value = 2478
$ws4.Range("F23").Value = 171
$ws4.Range("F25").Value = 567
$ws4.Range("F28").Value = 993
$ws4.Range("F29").Value = 1920
$ws4.Range("F30").Value = 1535
$ws4.Range("F32").Value = 1912
$ws4.Range("F35").Value = 20
$ws4.Range("F36").Value = 70
$ws4.Range("F39").Value = 52
$ws4.Range("F40").Value = 214
$ws4.Range("F41").Value = 389
$ws4.Range("F42").Value = 3
$ws4.Range("F46").Value = 53
$ws4.Range("F49").Value = 245
